# Copy pasted code. Now goes through both A AND B.
# Refactored the "build a sorted/merged sheet" logic into a reusable
# function instead of duplicating it: it walks a "primary" sheet's rows in
# order, and after each primary row, if the "secondary" sheet has a row
# whose column-C value matches, that row is written right underneath
# (leaving a blank separator row before the next primary record) - same
# layout "A Sorted" already used, just generalized to run for A->B and
# B->A.

function Build-SortedSheet {
    param($wb, $primaryWs, $secondaryWs, $newName, $afterWs)

    # Slurp the primary sheet's A/B/C columns (row 1..until blank) into memory.
    $primaryData = @()
    $r = 1
    while ($primaryWs.Cells.Item($r, 1).Value2 -ne $null) {
        $primaryData += , @($primaryWs.Cells.Item($r, 1).Value2, $primaryWs.Cells.Item($r, 2).Value2, $primaryWs.Cells.Item($r, 3).Value2)
        $r = $r + 1
    }

    # Index the secondary sheet by its column-C value so matches are a
    # simple lookup while we walk the primary rows.
    $secondaryLookup = @{}
    $r = 1
    while ($secondaryWs.Cells.Item($r, 1).Value2 -ne $null) {
        $cVal = $secondaryWs.Cells.Item($r, 3).Value2
        $secondaryLookup[$cVal] = @($secondaryWs.Cells.Item($r, 1).Value2, $secondaryWs.Cells.Item($r, 2).Value2, $secondaryWs.Cells.Item($r, 3).Value2)
        $r = $r + 1
    }

    $newWs = $wb.Worksheets.Add($null, $afterWs)
    $newWs.Name = $newName

    $outRow = 1
    foreach ($rec in $primaryData) {
        $newWs.Cells.Item($outRow, 1).Value = $rec[0]
        $newWs.Cells.Item($outRow, 2).Value = $rec[1]
        $newWs.Cells.Item($outRow, 3).Value = $rec[2]
        $outRow = $outRow + 1

        $cVal = $rec[2]
        if ($secondaryLookup.ContainsKey($cVal)) {
            $match = $secondaryLookup[$cVal]
            $newWs.Cells.Item($outRow, 1).Value = $match[0]
            $newWs.Cells.Item($outRow, 2).Value = $match[1]
            $newWs.Cells.Item($outRow, 3).Value = $match[2]
            $outRow = $outRow + 1
        }

        $outRow = $outRow + 1
    }

    return $newWs
}

$wb = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item("A")
$wsB = $wb.Worksheets.Item("B")
$wsASorted = $wb.Worksheets.Item("A Sorted")

# Old "A Sorted" behavior, preserved verbatim under its new name.
$wsASorted1 = Build-SortedSheet $wb $wsA $wsB "A Sorted1" $wsASorted

# New: the same merge, the other direction.
$wsBSorted = Build-SortedSheet $wb $wsB $wsA "B Sorted" $wsASorted1

Write-Output "Added sheets: A Sorted1, B Sorted"
